$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per the diff
$ws.Range("F3").Value = 5
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = 3
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -1
$ws.Range("F12").Value = -2
$ws.Range("F13").Value = -2
